$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rearranged/updated worker records (cols C=N° Doc, D=Nombre, E=Periodo Mora)
# Row 16 stays the same (CC / 45465928 / MERCEDES LOZANO GARCIA / 1709)
$ws.Range("C17").Value = "8785274"
$ws.Range("D17").Value = "MARCO HORACIO BELTRAN SERRANO"
$ws.Range("E17").Value = "1708"

$ws.Range("C18").Value = "1137223203"
$ws.Range("D18").Value = "JOSEFINA PERALES RAAD"
$ws.Range("E18").Value = "1708"

$ws.Range("C19").Value = "45465928"
$ws.Range("D19").Value = "MERCEDES LOZANO GARCIA"
$ws.Range("E19").Value = "1709"

$ws.Range("C20").Value = "8785274"
$ws.Range("D20").Value = "MARCO HORACIO BELTRAN SERRANO"
$ws.Range("E20").Value = "1709"

$ws.Range("C21").Value = "1137223203"
$ws.Range("D21").Value = "JOSEFINA PERALES RAAD"
$ws.Range("E21").Value = "1709"
